# Updated SO LUI TC with workaround for SOINV page opening issue
# Adds a new "Invoice" worksheet holding the InvoicePage_URL lookup pair,
# and makes it the active (selected) sheet.

$wb = $excel.ActiveWorkbook

# Add the new "Invoice" worksheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$invoiceSheet = $wb.Worksheets.Add($null, $lastSheet)
$invoiceSheet.Name = "Invoice"

# Populate the lookup values used to work around the SOINV page opening issue.
$invoiceSheet.Range("A1").Value = "InvoicePage_URL"
$invoiceSheet.Range("A2").Value = "https://rstk-dev-qa-ff.lightning.force.com/lightning/r/rstk__soinv__c/"

# Select cell F6 on the new sheet, matching the authored selection state.
$invoiceSheet.Range("F6").Select()

# Make the new Invoice sheet the active tab.
$invoiceSheet.Activate()
